# somani4_MP1.pptx edit:
#  1) Delete the trailing "Individual Contributions" slide (slide 13 / id 267).
#  2) Re-purpose the title slide: new title, trimmed author/affiliation subtitle.

$p = $ppt.ActivePresentation

# --- 1) Remove the last slide ("Individual Contributions") ---------------
$lastSlide = $p.Slides.Item($p.Slides.Count)
$lastSlide.Delete()

# --- 2) Update the title slide (slide 1) ----------------------------------
$titleSlide = $p.Slides.Item(1)

# Title placeholder: "Mini-Project 1 / ECE/CS 498DS / Spring 2020"
#   -> "Safety in Autonomous Vehicles"
$titleShape = $titleSlide.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Delete()
$titleShape.TextFrame.TextRange.Text = "Safety in Autonomous Vehicles"

# Subtitle placeholder: drop the team roster / credit-hours line, replace
# with the author name and affiliation.
$subtitleShape = $titleSlide.Shapes.Item(2)
$subtitleShape.TextFrame.TextRange.Text = "Akhilesh Somani (somani4)`rUniversity of Illinois at Urbana-Champaign"
